$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34; this shifts existing rows 34..159 down to 35..160
$ws.Rows(34).Insert()

# Populate the newly inserted row 34 with the new record
$ws.Range("A34").Value = 7
$ws.Range("B34").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C34").Value = 'Ñuble'
$ws.Range("D34").Value = 44459
$ws.Range("E34").Value = 16
$ws.Range("F34").Value = 100112023
$ws.Range("G34").Value = 'Brócoli'
$ws.Range("H34").Value = 'Sin especificar'
$ws.Range("I34").Value = 'Primera'
$ws.Range("J34").Value = 300
$ws.Range("K34").Value = 750
$ws.Range("L34").Value = 800
$ws.Range("M34").Value = 775
$ws.Range("N34").Value = '$/unidad'
$ws.Range("O34").Value = 'Provincia de Diguillín'
$ws.Range("P34").Value = 775
$ws.Range("Q34").Value = 1
$ws.Range("R34").Value = 'Hortaliza'
